# "Create api quote fee" — pretty-print the expected_validation_data JSON
# on the testcase sheet, which pushes row 2 taller and leaves the grid
# selection parked on a different cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testcase")

# F2 held a single-line JSON blob; reformat it to an indented / multi-line
# JSON literal (still one cell, wrapped across several physical lines).
$nl = [char]10
$json = '{' + $nl + '  "token": "NOT_NULL",' + $nl + '  "data.user_name": "phuongtt-chilinh"' + $nl + '}'
$ws.Range("F2").Value = $json

# The extra line breaks need more vertical room than the old 1-line text.
$ws.Rows.Item(2).RowHeight = 51.75

# Move the active selection (as last left by the editor) to G13.
$ws.Activate()
$ws.Range("G13").Select()
